$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.9682612609218907
$ws.Range("B2").Value = 0.7084336570791296
$ws.Range("C2").Value = 0.8958477337786327
$ws.Range("D2").Value = 6.108453273773193
$ws.Range("E2").Value = 6.585326235144985
$ws.Range("F2").Value = 6.108453273773193
$ws.Range("G2").Value = 0.07241352714325799
$ws.Range("H2").Value = -5.140192012851303
$ws.Range("I2").Value = 5.212605539994561
$ws.Range("K2").Value = 6.080389022827148
$ws.Range("L2").Value = 7
$ws.Range("P2").Value = 0.7139046192169189
$ws.Range("Q2").Value = 0.9527676701545715
$ws.Range("R2").Value = 0.1558623313903809
$ws.Range("S2").Value = 0.6651427149772644
$ws.Range("T2").Value = 0.962626039981842
$ws.Range("U2").Value = 0.1189597845077515
$ws.Range("V2").Value = 0.7326511740684509
$ws.Range("W2").Value = 0.9792417287826538
$ws.Range("X2").Value = 0.1084580421447754
$ws.Range("Y2").Value = 0.6869608163833618
$ws.Range("Z2").Value = 0.9856841564178467
$ws.Range("AA2").Value = 0.07648283243179321
$ws.Range("AB2").Value = 0.7651740312576294
$ws.Range("AC2").Value = 1.020087122917175
$ws.Range("AD2").Value = 0.03277677297592163
$ws.Range("AE2").Value = 0.6988017559051514
$ws.Range("AF2").Value = 1.008159995079041
$ws.Range("AG2").Value = 0.03411459922790527
$ws.Range("AH2").Value = 0.7353547811508179
$ws.Range("AI2").Value = 0.9822011590003967
$ws.Range("AJ2").Value = 0.1030478477478027
$ws.Range("AK2").Value = 0.6812342405319214
$ws.Range("AL2").Value = 0.9854488372802734
$ws.Range("AM2").Value = 0.0768769383430481
$ws.Range("AN2").Value = 0.7082946300506592
$ws.Range("AO2").Value = 0.9838261604309082
$ws.Range("AP2").Value = 0.09053057432174683
